$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 2 (this is a weekly price sheet; a new
# week's observation is prepended and everything else shifts down by one row).
$ws.Range("A2").EntireRow.Insert()

# The freshly inserted row picks up formatting copied from the row above
# (the bold/centered header style). Reset it to Normal so it matches the
# rest of the plain data rows, then re-apply the date number format that
# every other "Fecha" cell in column D uses.
$ws.Range("A2:T2").Style = "Normal"
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row with the latest market observation.
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "Vega Modelo de Temuco"
$ws.Range("C2").Value = "La Araucanía"
$ws.Range("D2").Value = 44699
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100108
$ws.Range("H2").Value = "Tropicales y subtropicales"
$ws.Range("I2").Value = 100108007
$ws.Range("J2").Value = "Coco"
$ws.Range("K2").Value = "Sin especificar"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 32000
$ws.Range("O2").Value = 32000
$ws.Range("P2").Value = 32000
$ws.Range("Q2").Value = "$/malla 20 unidades"
$ws.Range("R2").Value = "Perú"
$ws.Range("S2").Value = 1600
$ws.Range("T2").Value = 20
